$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct file endings: change file name extension from .xmi to .tsv
# in column A (the "file" column) for the data rows.
$ws.Range("A2").Value = "rwk1_digbib_300-1.tsv"
$ws.Range("A3").Value = "rwk1_digbib_1039-1.tsv"
$ws.Range("A4").Value = "rwk1_digbib_1057-1.tsv"
$ws.Range("A5").Value = "rwk1_mkhz_2778-1.tsv"
$ws.Range("A6").Value = "rwk1_mkhz_6147-1.tsv"
$ws.Range("A7").Value = "rwk1_mkhz_6263-1.tsv"

# Update the active sheet selection to match the saved workbook state
# (selection moved to column A, the file-name column that was corrected)
[void]$ws.Columns.Item(1).Select()
